$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6282199854184967
$ws.Range("C2").Value = 0.2173507237143042
$ws.Range("D2").Value = 0.09392477947358202
$ws.Range("F2").Value = 2.806569869732797
$ws.Range("G2").Value = 0.002553379781415717
$ws.Range("I2").Value = 1.427880147032816
$ws.Range("J2").Value = 0.2597227798853083
$ws.Range("K2").Value = 0.8378658212816106
$ws.Range("M2").Value = 0.4022742031768729
$ws.Range("N2").Value = 2.864061514756749
$ws.Range("B3").Value = 0.5944919923339569
$ws.Range("C3").Value = 0.2088828853370615
$ws.Range("D3").Value = 0.09216139151788383
$ws.Range("F3").Value = 2.797617076168777
$ws.Range("G3").Value = 0.002557186062454499
$ws.Range("I3").Value = 1.427821499067974
$ws.Range("J3").Value = 0.2578531681327547
$ws.Range("K3").Value = 0.7965394483557304
$ws.Range("M3").Value = 0.3901519390343822
$ws.Range("N3").Value = 2.880017338199572
$ws.Range("B4").Value = 0.5741370176503438
$ws.Range("C4").Value = 0.2038141067736348
$ws.Range("D4").Value = 0.0911230147626938
$ws.Range("F4").Value = 2.793359471321679
$ws.Range("G4").Value = 0.002559647669602131
$ws.Range("I4").Value = 1.428365173017546
$ws.Range("J4").Value = 0.2568267508066668
$ws.Range("K4").Value = 0.7716481230440593
$ws.Range("M4").Value = 0.3829343429189862
$ws.Range("N4").Value = 2.890566548421909
$ws.Range("B5").Value = 0.5659313498464655
$ws.Range("C5").Value = 0.2017812931034939
$ws.Range("D5").Value = 0.09071106263086648
$ws.Range("F5").Value = 2.791936063630473
$ws.Range("G5").Value = 0.002560682209019101
$ws.Range("I5").Value = 1.428732437600125
$ws.Range("J5").Value = 0.2564390376810479
$ws.Range("K5").Value = 0.7616262344729705
$ws.Range("M5").Value = 0.3800498560619801
$ws.Range("N5").Value = 2.895054509394306
$ws.Range("B6").Value = 0.5645741923147227
$ws.Range("C6").Value = 0.2014457224513393
$ws.Range("D6").Value = 0.09064333569096306
$ws.Range("F6").Value = 2.791718525948099
$ws.Range("G6").Value = 0.002560855893673195
$ws.Range("I6").Value = 1.428802220813481
$ws.Range("J6").Value = 0.2563765039973518
$ws.Range("K6").Value = 0.7599694481471602
$ws.Range("M6").Value = 0.3795743173129154
$ws.Range("N6").Value = 2.895811148746844
$ws.Range("B7").Value = 0.574025992052384
$ws.Range("C7").Value = 0.2037865589960433
$ws.Range("D7").Value = 0.09111741364767312
$ws.Range("F7").Value = 2.793339013185829
$ws.Range("G7").Value = 0.002559661494369367
$ws.Range("I7").Value = 1.428369536165576
$ws.Range("J7").Value = 0.2568213982333205
$ws.Range("K7").Value = 0.7715124722624864
$ws.Range("M7").Value = 0.3828952119291174
$ws.Range("N7").Value = 2.890626309221851
$ws.Range("B8").Value = 0.6165171142842212
$ws.Range("C8").Value = 0.2144038848487924
$ws.Range("D8").Value = 0.09330758298431618
$ws.Range("F8").Value = 2.803225647797134
$ws.Range("G8").Value = 0.002554666399363165
$ws.Range("I8").Value = 1.427739576110511
$ws.Range("J8").Value = 0.2590529125258598
$ws.Range("K8").Value = 0.8235161855120339
$ws.Range("M8").Value = 0.398047653458633
$ws.Range("N8").Value = 2.869406893363085
$ws.Range("B9").Value = 0.7026531003316165
$ws.Range("C9").Value = 0.236264172866953
$ws.Range("D9").Value = 0.09795286507952028
$ws.Range("F9").Value = 2.832455241537758
$ws.Range("G9").Value = 0.002545854648941435
$ws.Range("I9").Value = 1.431107684036654
$ws.Range("J9").Value = 0.264393763386515
$ws.Range("K9").Value = 0.9293354592611252
$ws.Range("M9").Value = 0.429551676701351
$ws.Range("N9").Value = 2.833768770551004
$ws.Range("B10").Value = 0.7676599757377289
$ws.Range("C10").Value = 0.252966630831736
$ws.Range("D10").Value = 0.1015777743868398
$ws.Range("F10").Value = 2.859946865173498
$ws.Range("G10").Value = 0.002539973966593725
$ws.Range("I10").Value = 1.436396263493876
$ws.Range("J10").Value = 0.2689074739910922
$ws.Range("K10").Value = 1.009439611950853
$ws.Range("M10").Value = 0.4537931142658991
$ws.Range("N10").Value = 2.811232565332418
$ws.Range("B11").Value = 0.797609838250537
$ws.Range("C11").Value = 0.2607061799464816
$ws.Range("D11").Value = 0.1032725800959611
$ws.Range("F11").Value = 2.873764081774809
$ws.Range("G11").Value = 0.002537426175963868
$ws.Range("I11").Value = 1.439415049255366
$ws.Range("J11").Value = 0.2710893462819541
$ws.Range("K11").Value = 1.046397682324169
$ws.Range("M11").Value = 0.4650601367326246
$ws.Range("N11").Value = 2.801773537998471
$ws.Range("B12").Value = 0.8090054539072185
$ws.Range("C12").Value = 0.26365740315984
$ws.Range("D12").Value = 0.1039209126181362
$ws.Range("F12").Value = 2.879185067255563
$ws.Range("G12").Value = 0.002536479606212155
$ws.Range("I12").Value = 1.440646446833604
$ws.Range("J12").Value = 0.2719340694054182
$ws.Range("K12").Value = 1.060467433238756
$ws.Range("M12").Value = 0.4693611263014006
$ws.Range("N12").Value = 2.798305789047902
$ws.Range("B13").Value = 0.8065487926808714
$ws.Range("C13").Value = 0.2630208948868358
$ws.Range("D13").Value = 0.1037809921219264
$ws.Range("F13").Value = 2.878009165620711
$ws.Range("G13").Value = 0.002536682657912615
$ws.Range("I13").Value = 1.440377316776917
$ws.Range("J13").Value = 0.2717513204396624
$ws.Range("K13").Value = 1.057433945739518
$ws.Range("M13").Value = 0.4684333010800401
$ws.Range("N13").Value = 2.799047549835564
$ws.Range("B14").Value = 0.7985462756091977
$ws.Range("C14").Value = 0.2609485691020268
$ws.Range("D14").Value = 0.1033257878429339
$ws.Range("F14").Value = 2.874206286609464
$ws.Range("G14").Value = 0.002537347936380046
$ws.Range("I14").Value = 1.439514587839568
$ws.Range("J14").Value = 0.27115847143601
$ws.Range("K14").Value = 1.047553715444622
$ws.Range("M14").Value = 0.4654132921778
$ws.Range("N14").Value = 2.801485954863139
$ws.Range("B15").Value = 0.7936515655153755
$ws.Range("C15").Value = 0.259681871737456
$ws.Range("D15").Value = 0.1030478133705088
$ws.Range("F15").Value = 2.871901493423138
$ws.Range("G15").Value = 0.002537757809196356
$ws.Range("I15").Value = 1.438997637154372
$ws.Range("J15").Value = 0.2707977432109203
$ws.Range("K15").Value = 1.041511495475675
$ws.Range("M15").Value = 0.4635679304151523
$ws.Range("N15").Value = 2.802994423143403
$ws.Range("B16").Value = 0.7657102785103973
$ws.Range("C16").Value = 0.2524636870749646
$ws.Range("D16").Value = 0.1014679331963322
$ws.Range("F16").Value = 2.85907027504075
$ws.Range("G16").Value = 0.002540143025869134
$ws.Range("I16").Value = 1.436211322264143
$ws.Range("J16").Value = 0.2687674712115324
$ws.Range("K16").Value = 1.00703474455716
$ws.Range("M16").Value = 0.4530616046245228
$ws.Range("N16").Value = 2.811866703144389
$ws.Range("B17").Value = 0.7486659304512671
$ws.Range("C17").Value = 0.2480718627530223
$ws.Range("D17").Value = 0.100510432584926
$ws.Range("F17").Value = 2.851534668170501
$ws.Range("G17").Value = 0.002541638833834255
$ws.Range("I17").Value = 1.434659080503174
$ws.Range("J17").Value = 0.2675548988092231
$ws.Range("K17").Value = 0.9860171213959461
$ws.Range("M17").Value = 0.4466776315367653
$ws.Range("N17").Value = 2.817512739205853
$ws.Range("B18").Value = 0.7388980480288296
$ws.Range("C18").Value = 0.2455591173385017
$ws.Range("D18").Value = 0.0999640184929973
$ws.Range("F18").Value = 2.847323786767873
$ws.Range("G18").Value = 0.002542511176530262
$ws.Range("I18").Value = 1.433823962786001
$ws.Range("J18").Value = 0.2668695596374704
$ws.Range("K18").Value = 0.9739771239302115
$ws.Range("M18").Value = 0.4430282804904451
$ws.Range("N18").Value = 2.820834801641453
$ws.Range("B19").Value = 0.7355969238433602
$ws.Range("C19").Value = 0.2447106296805543
$ws.Range("D19").Value = 0.09977975444664366
$ws.Range("F19").Value = 2.845919243878527
$ws.Range("G19").Value = 0.002542808599459392
$ws.Range("I19").Value = 1.433551111633278
$ws.Range("J19").Value = 0.2666395935818002
$ws.Range("K19").Value = 0.9699089662667859
$ws.Range("M19").Value = 0.4417965451636263
$ws.Range("N19").Value = 2.821972404309307
$ws.Range("B20").Value = 0.7504766491019836
$ws.Range("C20").Value = 0.2485380012857092
$ws.Range("D20").Value = 0.1006119138858139
$ws.Range("F20").Value = 2.852324073722258
$ws.Range("G20").Value = 0.002541478362048228
$ws.Range("I20").Value = 1.434818347777593
$ws.Range("J20").Value = 0.2676827268075073
$ws.Range("K20").Value = 0.9882494343504789
$ws.Range("M20").Value = 0.4473548835887726
$ws.Range("N20").Value = 2.816903985716735
$ws.Range("B21").Value = 0.8008953371580958
$ws.Range("C21").Value = 0.2615567065240612
$ws.Range("D21").Value = 0.1034593150391032
$ws.Range("F21").Value = 2.875318161602124
$ws.Range("G21").Value = 0.002537152034164125
$ws.Range("I21").Value = 1.439765596360886
$ws.Range("J21").Value = 0.2713321035297156
$ws.Range("K21").Value = 1.050453755773731
$ws.Range("M21").Value = 0.4662994082239038
$ws.Range("N21").Value = 2.800766635986008
$ws.Range("B22").Value = 0.8341630510235518
$ws.Range("C22").Value = 0.2701842644858345
$ws.Range("D22").Value = 0.1053583943916578
$ws.Range("F22").Value = 2.891446127582924
$ws.Range("G22").Value = 0.002534430705496373
$ws.Range("I22").Value = 1.443513340960038
$ws.Range("J22").Value = 0.2738249919983389
$ws.Range("K22").Value = 1.091542307780173
$ws.Range("M22").Value = 0.4788813240577028
$ws.Range("N22").Value = 2.790885554299351
$ws.Range("B23").Value = 0.8163785439610649
$ws.Range("C23").Value = 0.2655686537167981
$ws.Range("D23").Value = 0.1043413441360599
$ws.Range("F23").Value = 2.882737622482381
$ws.Range("G23").Value = 0.002535873444549963
$ws.Range("I23").Value = 1.441465994086109
$ws.Range("J23").Value = 0.2724846230196221
$ws.Range("K23").Value = 1.069572821919934
$ws.Range("M23").Value = 0.4721477734770332
$ws.Range("N23").Value = 2.796098311376753
$ws.Range("B24").Value = 0.7496579261944021
$ws.Range("C24").Value = 0.24832722214407
$ws.Range("D24").Value = 0.1005660215229938
$ws.Range("F24").Value = 2.851966805191537
$ws.Range("G24").Value = 0.002541550872873509
$ws.Range("I24").Value = 1.434746164594529
$ws.Range("J24").Value = 0.2676248990622554
$ws.Range("K24").Value = 0.9872400707462248
$ws.Range("M24").Value = 0.4470486328961201
$ws.Range("N24").Value = 2.817178966168939
$ws.Range("B25").Value = 0.6790491656012421
$ws.Range("C25").Value = 0.230238267157489
$ws.Range("D25").Value = 0.09665884778713973
$ws.Range("F25").Value = 2.823492506263449
$ws.Range("G25").Value = 0.002548133814659872
$ws.Range("I25").Value = 1.429702909587959
$ws.Range("J25").Value = 0.2628454652871781
$ws.Range("K25").Value = 0.9002956152561978
$ws.Range("M25").Value = 0.4208369372125773
$ws.Range("N25").Value = 2.84276978667279
